# Highlight the Matrix2D.c implementation bullets on "Step 2" (slide 9).
#
# Four sub-bullets get a green highlight and the "test every function"
# bullet gets a yellow highlight, matching the completed-implementation
# checklist called out in the commit message.

function RGB($r, $g, $b) {
    return $r + ($g * 256) + ($b * 65536)
}

$Green  = RGB 0 255 0
$Yellow = RGB 255 255 0

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(9)
$shp = $s.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange
$fullText = $tr.Text

$targets = @(
    @{ Text = "Implement the Matrix2DIdentity function"; Color = $Green },
    @{ Text = "Implement the Matrix2DTranslate, Scale, and Rot* functions"; Color = $Green },
    @{ Text = "Implement the Matrix2DConcat function"; Color = $Green },
    @{ Text = "Implement the remaining functions at any time"; Color = $Green },
    @{ Text = "Make sure to test every single function!"; Color = $Yellow }
)

foreach ($target in $targets) {
    $idx = $fullText.IndexOf($target.Text)
    $run = $tr.Characters($idx + 1, $target.Text.Length)
    $run.Font.Highlight = $target.Color
}
